$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "2025-03-13 00:00:00",
    "2025-04-07 00:00:00",
    "2025-05-05 00:00:00",
    "2025-04-30 00:00:00",
    "2025-05-08 00:00:00",
    "2025-04-14 00:00:00",
    "2025-04-24 00:00:00",
    "2025-05-16 00:00:00",
    "2025-04-28 00:00:00",
    "2025-05-02 00:00:00",
    "2025-04-28 00:00:00",
    "2025-04-15 00:00:00",
    "2025-04-23 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-06 00:00:00",
    "2025-05-06 00:00:00",
    "2025-05-15 00:00:00",
    "2025-04-23 00:00:00",
    "2025-05-09 00:00:00",
    "2025-04-22 00:00:00",
    "2025-05-15 00:00:00",
    "2025-05-12 00:00:00",
    "2025-04-16 00:00:00",
    "2025-04-23 00:00:00",
    "2025-04-16 00:00:00",
    "2025-05-14 00:00:00",
    "2025-05-05 00:00:00",
    "2025-05-12 00:00:00",
    "2025-04-23 00:00:00",
    "2025-05-07 00:00:00",
    "2025-05-07 00:00:00",
    "2025-05-07 00:00:00",
    "2025-04-30 00:00:00",
    "2025-04-30 00:00:00",
    "2025-05-05 00:00:00",
    "2025-05-30 00:00:00",
    "2025-06-10 00:00:00",
    "2025-12-31 00:00:00",
    "2025-04-28 00:00:00",
    "2025-12-31 00:00:00",
    "2023-11-06 00:00:00",
    "2025-04-28 00:00:00",
    "2025-04-28 00:00:00",
    "2025-04-28 00:00:00",
    "2025-04-30 00:00:00",
    "2025-05-07 00:00:00",
    "2025-05-08 00:00:00",
    "2025-05-08 00:00:00",
    "2025-05-08 00:00:00",
    "2025-05-16 00:00:00",
    "2025-05-22 00:00:00",
    "2024-09-16 00:00:00",
    "2025-04-18 00:00:00",
    "2025-01-25 00:00:00",
    "2025-04-23 00:00:00",
    "2025-04-28 00:00:00",
    "2025-04-28 00:00:00",
    "2025-05-05 00:00:00",
    "2025-05-06 00:00:00",
    "2024-08-02 00:00:00",
    "2025-05-19 00:00:00",
    "2025-04-18 00:00:00",
    "2024-09-30 00:00:00",
    "2025-03-18 00:00:00",
    "2025-03-28 00:00:00",
    "2025-05-05 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-09 00:00:00",
    "2025-05-15 00:00:00",
    "2025-05-15 00:00:00",
    "2025-05-15 00:00:00",
    "2025-05-15 00:00:00",
    "2025-05-12 00:00:00",
    "2025-05-22 00:00:00",
    "2025-05-23 00:00:00",
    "2025-05-26 00:00:00",
    "2025-05-23 00:00:00",
    "2025-05-23 00:00:00",
    "2025-06-06 00:00:00",
    "2025-09-01 00:00:00",
    "2025-06-19 00:00:00",
    "2025-12-31 00:00:00",
    "2025-12-31 00:00:00",
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $values[$i]
}